$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 02:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1591464
$ws.Range("C4").Value = 20881
$ws.Range("D4").Value = 370068
$ws.Range("E4").Value = 1126462
$ws.Range("G4").Value = 1401
$ws.Range("H4").Value = 94934

# Brasil (row 6)
$ws.Range("B6").Value = 293357
$ws.Range("C6").Value = 21472
$ws.Range("E6").Value = 157780
$ws.Range("G6").Value = 911
$ws.Range("H6").Value = 18894

# Canada (row 17)
$ws.Range("B17").Value = 80142
$ws.Range("C17").Value = 1030
$ws.Range("D17").Value = 40776
$ws.Range("E17").Value = 33335
$ws.Range("G17").Value = 119
$ws.Range("H17").Value = 6031

# Ghana (row 64)
$ws.Range("B64").Value = 6269
$ws.Range("C64").Value = 173
$ws.Range("D64").Value = 1898
$ws.Range("E64").Value = 4340

# Republica de Africa Central overtakes Estado de Palestina, Etiopia and
# Madagascar in the case-count ranking, so those four countries shift down
# one row (rows 135-138 keep their position, the country label + stats move).
$ws.Range("A135").Value = "Republica de Africa Central"
$ws.Range("B135").Value = 418
$ws.Range("C135").Value = 52
$ws.Range("D135").Value = 18
$ws.Range("E135").Value = 400
$ws.Range("H135").Value = 0

$ws.Range("A136").Value = "Estado de Palestina"
$ws.Range("B136").Value = 398
$ws.Range("C136").Value = 7
$ws.Range("D136").Value = 346
$ws.Range("E136").Value = 50
$ws.Range("H136").Value = 2

$ws.Range("A137").Value = "Etiopia"
$ws.Range("B137").Value = 389
$ws.Range("C137").Value = 24
$ws.Range("D137").Value = 122
$ws.Range("E137").Value = 262
$ws.Range("H137").Value = 5

$ws.Range("A138").Value = "Madagascar"
$ws.Range("B138").Value = 371
$ws.Range("C138").Value = 45
$ws.Range("D138").Value = 131
$ws.Range("E138").Value = 238
$ws.Range("H138").Value = 2

# Bahamas ties/overtakes Monaco, so they swap rows (row 170 <-> row 171).
$ws.Range("A170").Value = "Bahamas"
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 43
$ws.Range("E170").Value = 43
$ws.Range("H170").Value = 11

$ws.Range("A171").Value = "Monaco"
$ws.Range("B171").Value = 97
$ws.Range("D171").Value = 90
$ws.Range("E171").Value = 3
$ws.Range("H171").Value = 4
